# ============================================================================
# Edit: insert a "2022-Q1" sheet (fund-holding detail) before the "总计"
# sheet. The previous "总计" worksheet's content becomes the new "2022-Q1"
# sheet (same sheetId/relationship slot), and a fresh "总计" worksheet is
# appended right after it with the historical roll-up table plus a new
# leading row for "2022-Q1".
# ============================================================================

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Formatting donors: plain cells copied later to (re)stamp formats without
# fabricating brand-new style entries.
#   $plainDonor  -> default (un-styled) cell, used to clear the temporary
#                   Text ("@") format once literal numeric-looking strings
#                   have been written (fund codes, percentages, ...).
#   $boldDonor   -> the bold/bordered index-column style (column A cells).
#   $headerDonor -> the bold/bordered header style (row 1 cells).
# ---------------------------------------------------------------------------
$refSheet    = $wb.Worksheets.Item("2021-Q4")
$plainDonor  = $refSheet.Range("B2")
$boldDonor   = $refSheet.Range("A2")
$headerDonor = $refSheet.Range("B1")

# ---------------------------------------------------------------------------
# Step 1: the existing "总计" sheet keeps its sheetId/relationship slot but
# becomes "2022-Q1" with a brand-new fund-holdings table.
# ---------------------------------------------------------------------------
$q1 = $wb.Worksheets.Item("总计")
$q1.Cells.Clear()
$q1.Name = "2022-Q1"

# Header row (B1:H1)
$q1.Range("B1").Value = "基金代码"
$q1.Range("C1").Value = "基金名称"
$q1.Range("D1").Value = "基金规模"
$q1.Range("E1").Value = "股票总仓位"
$q1.Range("F1").Value = "仓位占比"
$q1.Range("G1").Value = "持有市值(亿元)"
$q1.Range("H1").Value = "仓位排名"

# Data rows 2..34 -- column A (running index) and column H (rank) are
# numeric; columns B,D,E,F,G are numeric-looking strings that must stay text
# (fund code, scale, position, ratio, market value), so that block is marked
# Text ("@") before writing and restored to the default format afterwards.
# Column C (fund name) is plain (non-numeric) text throughout.
$idxArr   = @("0","1","2","3","4","5","6","7","8","9","10","11","12","13","14","15","16","17","18","19","20","21","22","23","24","25","26","27","28","29","30","31","32")
$codeArr  = @("501011","011882","012123","010996","012196","001869","005505","501012","005506","005689","161611","012197","004569","161123","161729","007718","003581","515950","010997","011288","011883","400013","007687","159838","501028","007111","010159","011289","585001","014157","010500","007613","005443")
$nameArr  = @("汇添富中证中药指数（LOF）A","招商蓝筹精选股票型证券投资基金A","招商金安成长严选1年封闭运作混合型证券投资基金","招商品质升级混合A","招商品质生活混合型证券投资基金A","招商制造业转型灵活配置混合A","前海开源中药研究精选股票A","汇添富中证中药指数（LOF）C","前海开源中药研究精选股票C","中银医疗保健灵活配置混合A","融通内需驱动混合","招商品质生活混合型证券投资基金C","招商制造业转型灵活配置混合C","易方达并购重组指数（LOF）","招商 3 年封闭运作瑞利灵活配置混合型","中银创新医疗混合A","新疆前海联合国民健康产业灵活配置混合A","富国中证医药50ETF","招商品质升级混合C","上银医疗健康混合A","招商蓝筹精选股票型证券投资基金C","东方成长收益灵活配置混合A","东方成长收益灵活配置混合C","博时中证医药50交易型开放式指数证券投资基金","财通多策略福瑞混合（LOF）","新疆前海联合国民健康产业灵活配置混合C","中银医疗保健灵活配置混合C","上银医疗健康混合C","东吴中证新兴","国泰君安创新医药混合","中银创新医疗混合C","嘉合医疗健康混合","国金量化多策略灵活配置混合")
$scaleArr = @("24.13","48.27","36.13","25.62","25.47","26.05","11.44","8.91","4.62","7.30","12.96","7.73","6.38","4.78","5.33","3.59","2.66","4.03","2.76","1.48","1.34","2.43","2.15","1.26","2.14","0.32","0.32","0.23","0.56","0.31","0.13","0.29","0.51")
$posArr   = @("94.57","86.18","86.22","85.13","86.47","86.34","91.95","94.57","91.95","90.43","65.58","86.47","86.34","94.71","86.78","89.91","94.65","99.22","85.13","81.70","86.18","53.15","53.15","98.62","55.41","94.65","90.43","81.70","93.03","79.44","89.91","92.44","64.10")
$ratioArr = @("7.42","3.64","3.35","3.68","3.61","3.37","7.48","7.42","7.48","4.48","2.50","3.61","3.37","4.47","3.42","5.02","4.62","2.86","3.68","5.10","3.64","1.70","1.70","2.85","1.53","4.62","4.48","5.10","1.71","2.88","5.02","1.90","0.94")
$mvArr    = @("1.7904","1.7570","1.2104","0.9428","0.9195","0.8779","0.8557","0.6611","0.3456","0.3270","0.3240","0.2791","0.2150","0.2137","0.1823","0.1802","0.1229","0.1153","0.1016","0.0755","0.0488","0.0413","0.0366","0.0359","0.0327","0.0148","0.0143","0.0117","0.0096","0.0089","0.0065","0.0055","0.0048")
$rankArr  = @("1","8","9","9","8","10","8","1","8","3","6","8","10","6","9","3","10","10","9","3","8","8","8","10","9","10","3","3","2","8","3","10","1")

$dataRange = $q1.Range("B2:G34")
$dataRange.NumberFormat = "@"

for ($i = 0; $i -lt 33; $i++) {
    $r = 2 + $i
    $q1.Range("A$r").Value = [double]$idxArr[$i]
    $q1.Range("B$r").Value = $codeArr[$i]
    $q1.Range("C$r").Value = $nameArr[$i]
    $q1.Range("D$r").Value = $scaleArr[$i]
    $q1.Range("E$r").Value = $posArr[$i]
    $q1.Range("F$r").Value = $ratioArr[$i]
    $q1.Range("G$r").Value = $mvArr[$i]
    $q1.Range("H$r").Value = [double]$rankArr[$i]
}

# Restore the default (un-styled) format on the text block we wrote.
$plainDonor.Copy()
$dataRange.PasteSpecial(-4122)

# Header row (B1:H1) + index column (A2:A34) use the bold/bordered style.
$headerDonor.Copy()
$q1.Range("B1:H1").PasteSpecial(-4122)
$boldDonor.Copy()
$q1.Range("A2:A34").PasteSpecial(-4122)

$q1.Range("A1").Select()
$excel.CutCopyMode = $false

# ---------------------------------------------------------------------------
# Step 2: append a fresh "总计" sheet right after "2022-Q1" holding the
# historical roll-up table (same columns as before) plus the new "2022-Q1"
# row at the top of the data.
# ---------------------------------------------------------------------------
$total = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $q1)
$total.Name = "总计"

# Header row (B1:D1)
$total.Range("B1").Value = "日期"
$total.Range("C1").Value = "持有数量(只)"
$total.Range("D1").Value = "持有市值(亿元)"

$tIdxArr   = @("0","1","2","3","4","5")
$tLabelArr = @("2022-Q1","2021-Q4","2021-Q3","2021-Q2","2021-Q1","2020-Q4")
$tCountArr = @("33","64","37","45","60","72")
$tMvArr    = @("11.77","17.58","9.39","5.59","18.08","10.52")

# Column B (quarter label, e.g. "2022-Q1") stays plain text naturally (it is
# not numeric-looking, Excel will not coerce it), columns A/C/D are numbers.
for ($i = 0; $i -lt 6; $i++) {
    $r = 2 + $i
    $total.Range("A$r").Value = [double]$tIdxArr[$i]
    $total.Range("B$r").Value = $tLabelArr[$i]
    $total.Range("C$r").Value = [double]$tCountArr[$i]
    $total.Range("D$r").Value = [double]$tMvArr[$i]
}

# Header row (B1:D1) + index column (A2:A7) use the bold/bordered style.
$headerDonor.Copy()
$total.Range("B1:D1").PasteSpecial(-4122)
$boldDonor.Copy()
$total.Range("A2:A7").PasteSpecial(-4122)

$total.Range("A1").Select()
$excel.CutCopyMode = $false
